$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New variable label used across all data rows
$newVar = "Diferença 2024/03 - 2023/03"

# Row 2: Roraima (unchanged region/rank, new value + label)
$ws.Range("B2").Value = $newVar
$ws.Range("C2").Value = 3.179824561403514

# Row 3: Mato Grosso (unchanged region/rank, new value + label)
$ws.Range("B3").Value = $newVar
$ws.Range("C3").Value = 2.762760003629232

# Row 4: now Paraíba (was Rio Grande do Norte), rank stays 3º
$ws.Range("A4").Value = "Paraíba"
$ws.Range("B4").Value = $newVar
$ws.Range("C4").Value = 2.662615874305786

# Row 5: now Rio Grande do Norte (was Acre), rank stays 4º
$ws.Range("A5").Value = "Rio Grande do Norte"
$ws.Range("B5").Value = $newVar
$ws.Range("C5").Value = 2.490310519356157

# Row 6: now Acre (was Sergipe), rank stays 5º
$ws.Range("A6").Value = "Acre"
$ws.Range("B6").Value = $newVar
$ws.Range("C6").Value = 2.323730701893993

# Row 7: now Tocantins (was Paraíba), rank stays 6º
$ws.Range("A7").Value = "Tocantins"
$ws.Range("B7").Value = $newVar
$ws.Range("C7").Value = 2.264465678125013

# Row 8: now Sergipe (was Nordeste), new value, gains rank 8º
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = $newVar
$ws.Range("C8").Value = 1.285889996082624
$ws.Range("D8").Value = "8º"

# Row 9: new row - Nordeste (shifted down from old row 8, no rank)
$ws.Range("A9").Value = "Nordeste"
$ws.Range("B9").Value = $newVar
$ws.Range("C9").Value = 0.4591818090789772

# Row 10: new row - Brasil (shifted down from old row 9, no rank)
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = $newVar
$ws.Range("C10").Value = 0.5397773419660581
